$d = $word.ActiveDocument

function Replace-Text($find, $replace) {
    $d.Content.Find.Execute($find, $true, $false, $false, $false, $false,
                             $true, 1, $false, $replace, 2)
}

Replace-Text "media quieries" "media queries"
Replace-Text "position:absolute;" "position: absolute;"
Replace-Text "W3 svg arrows" "W3 SVG arrows"
Replace-Text "from wordpress and deplaying" "from Wordpress and deploying"
Replace-Text "with bluehost ghosting" "with Bluehost ghosting"
Replace-Text "was getting wordpress to be shown" "was getting Wordpress to be shown"
Replace-Text "fetched on netlify. Really furstrating." "fetched on Netlify. Really frustrating."
